$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix: Redirected Precision-Recall curve to (1, min(Precision) = P/(N+P))
# -> apply a custom number format to E4 (the min-Precision value cell)
$ws.Range("E4").NumberFormat = "#,##0.0000_);\(#,##0.0000\)"

# Updated performance numbers resulting from the fix
$ws.Range("F7").Value = 0.90369999999999995
$ws.Range("H7").Value = 0.9486
$ws.Range("D13").Value = 0.94159999999999999

# Removed the "lower/upper bound of CI for optimal F1 measure" rows
# (rows 14 & 15), which shifts the "weighted sum / binary-classifier score"
# row up from row 16 to row 14.
$ws.Rows("14:15").Delete()

# Reflect where the user ended up after editing: scrolled down a bit and
# selected the last data cell.
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("H14").Select() | Out-Null
